# Edit: Randomly allocate people with flexible time preference ("w") to a
# specific week ("m") group rather than leaving them discriminated.
# Concretely: every cell in column D (the week-availability code) whose
# value is "w" (the flexible/"wahlweise" marker) is re-labeled "nm" so
# people with specific week availability are no longer mixed in with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds the availability code for each person (rows 2-12).
# Replace every "w" with "nm" - leave "m" entries untouched.
for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    if ($cell.Value2 -eq "w") {
        $cell.Value2 = "nm"
    }
}

# Update the active cell selection to A13 (matches the saved view state).
$ws.Range("A13").Select()
